# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F9"  = 97372
        "F10" = 813
        "F11" = 40492
        "F14" = 616
        "F16" = 722
        "F17" = 1373
        "F18" = 231
        "F19" = 442
        "F22" = 5416
        "F25" = 2524
        "F26" = 6014
        "F28" = 1064
        "F29" = 634
        "F37" = 761
        "F39" = 65
        "F44" = 79
        "F47" = 701
    }
    "演出" = @{
        "F7"  = 696
        "F24" = 157
        "F31" = 102
        "F32" = 246
        "F36" = 46
        "F38" = 5
    }
    "本地生活" = @{
        "F6" = 465
    }
    "全部类型" = @{
        "F6"  = 465
        "F13" = 813
        "F14" = 696
        "F15" = 1373
        "F16" = 231
        "F21" = 2524
        "F22" = 6014
        "F24" = 1064
        "F27" = 634
        "F33" = 761
        "F36" = 65
        "F40" = 102
        "F43" = 46
        "F45" = 5
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $updates[$sheetName]
    foreach ($cellRef in $cells.Keys) {
        $ws.Range($cellRef).Value = $cells[$cellRef]
    }
}
